# Update EDCR Results sheet with new nickel results.
# Rows 2-100 get uniform new values in columns B, C, D, F, H.
# Columns A (epsilon), E (NSC), G (NRC) are left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EDCR Results")

$lastRow = 100

$ws.Range("B2:B$lastRow").Value = 0.18125
$ws.Range("C2:C$lastRow").Value = 0.4603174603174603
$ws.Range("D2:D$lastRow").Value = 0.2600896860986547
$ws.Range("F2:F$lastRow").Value = 44
$ws.Range("H2:H$lastRow").Value = 8
